$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") from 45190 to 45192 for all data rows (2-533)
$ws.Range("C2:C533").Value = 45192
